# Update p-value table ("Cod: p-values for trend lines in Fig 8") in the
# "ecological" column: replace old p-values with the new ones from the
# commit "no OVER, yes RECOVERY".
#
# Mapping (row -> old -> new):
#   GDP 2016      : 0.07 -> 0.43
#   OHI economic  : 0.92 -> 0.68
#   OHI fisheries : 0.88 -> 0.36
#   Readiness     : 1.00 -> 0.65
#   Vulnerability : 0.96 -> 0.45

$d = $word.ActiveDocument

function Replace-ExactText($old, $new) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-ExactText "0.07" "0.43"
Replace-ExactText "0.92" "0.68"
Replace-ExactText "0.88" "0.36"
Replace-ExactText "1.00" "0.65"
Replace-ExactText "0.96" "0.45"
